$d = $word.ActiveDocument

# Locate the paragraph that ends with the "Choose ONE..." sentence so the new
# citation paragraph can be inserted immediately after it. Keep reusing this
# same Range object -- Find.Execute mutates it in place to the found text,
# and re-fetching $d.Content afterwards would reset us back to the start of
# the document.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Choose ONE. All name pairs are from validated hiring discrimination research.",
                                    $true, $false, $false, $false, $false,
                                    $true, 1, $false, "", 0)

$srcPara  = $searchRange.Paragraphs(1)
$srcRange = $srcPara.Range

# Insert a brand-new (empty) paragraph right after it.
$srcRange.InsertParagraphAfter() | Out-Null
$citationPara  = $srcPara.Next()
$citationRange = $citationPara.Range

$quote1 = [char]0x201C
$quote2 = [char]0x201D

$runsXml =
  '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Citation:</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">Crabtree, Charles et al. </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">' + $quote1 + 'Validated names for experimental studies on race and ethnicity.' + $quote2 + '</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Scientific Data</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">vol. 10, no. 1 (2023): 130. https://doi.org/10.1038/s41597-023-01947-0</w:t></w:r>'

$xmlFrag =
  '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' + $runsXml + '</w:p></w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$citationRange.InsertXML($xmlFrag) | Out-Null

Write-Output ("Citation paragraph text: [" + $citationPara.Range.Text + "]")
Write-Output ("Citation paragraph style: " + $citationPara.Range.ParagraphStyle.NameLocal)
